# Applies the updated crypto price/volume snapshot values to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.341.55"
$ws.Range("E2").Value = "  -2.91%  "
$ws.Range("D3").Value = "3.346.06"
$ws.Range("E3").Value = "  -4.74%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'554.43"
$ws.Range("E5").Value = "  -5.07%  "
$ws.Range("D6").Value = "'175.02"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("D7").Value = "'0.616"
$ws.Range("E7").Value = "  -3.37%  "
$ws.Range("D8").Value = "3.339.62"
$ws.Range("E8").Value = "  -4.55%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "'0.626"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").Value = "'0.162"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "'54.57"
$ws.Range("E12").Value = "  -2.63%  "
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("D14").Value = "'9.03"
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("D15").Value = "3.894.11"
$ws.Range("E15").Value = "  -4.43%  "
$ws.Range("D16").Value = "'18.27"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.118"
$ws.Range("E17").Value = "  -3.15%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.356.83"
$ws.Range("E18").Value = "  -4.62%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "64.364.29"
$ws.Range("E19").Value = "  -2.84%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'11.76"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").Value = "'0.978"
$ws.Range("D22").Value = "'435.34"
$ws.Range("E22").Value = "  +4.77%  "
$ws.Range("E23").Value = "  +12.11%  "
$ws.Range("D24").Value = "'4.07"
$ws.Range("E24").Value = "  -5.19%  "
$ws.Range("D25").Value = "'84.34"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").Value = "'13.41"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "'10.73"
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("D28").Value = "'2.82"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("D29").Value = "'8.70"
$ws.Range("E29").Value = "  -5.14%  "
$ws.Range("D30").Value = "'29.63"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").Value = "'11.44"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("D33").Value = "'576.11"
$ws.Range("E33").Value = "  -3.19%  "
$ws.Range("D34").Value = "'0.107"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("D35").Value = "'58.31"
$ws.Range("E35").Value = "  -4.16%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  -8.63%  "
$ws.Range("D38").Value = "'3.52"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("D39").Value = "'35.58"
$ws.Range("E39").Value = "  -3.65%  "
$ws.Range("D40").Value = "0.0₃0751"
$ws.Range("E40").Value = "  -5.64%  "
$ws.Range("D41").Value = "'0.366"
$ws.Range("E41").Value = "  -5.35%  "
$ws.Range("D42").Value = "3.097.78"
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "'2.80"
$ws.Range("E44").Value = "  -5.90%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0409"
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.21"
$ws.Range("E46").Value = "  -4.52%  "
$ws.Range("D47").Value = "'2.45"
$ws.Range("E47").Value = "  -3.85%  "
$ws.Range("E48").Value = "  -2.77%  "
$ws.Range("D49").Value = "'2.59"
$ws.Range("E49").Value = "  -3.46%  "
$ws.Range("D50").Value = "'135.42"
$ws.Range("E50").Value = "  -3.42%  "
$ws.Range("D51").Value = "'8.27"
$ws.Range("E51").Value = "  -4.03%  "
